# Bug fix to reg_RSME.xlsx: remove erroneous duplicated IT rows that were
# mistakenly included at the top of the UK sheet (rows 2-5), shifting the
# remaining UK rows up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UK")

# Delete rows 2 through 5 (inclusive) - these contained incorrect IT entries.
$ws.Rows("2:5").Delete()
